# Remove the obsolete footnote markers (e.g. " [5]") from vaccine/brand labels
# and collapse the multi-line cell labels (which used embedded line breaks) down
# to single-line text, per the source-data cleanup described in the commit message.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Pediatric VFC Vaccine ')
$ws.Range('A2').Value = 'Dengue Tetravalent Vaccine, Live '
$ws.Range('A3').Value = 'DTaP '
$ws.Range('A4').Value = 'DTaP '
$ws.Range('A5').Value = 'DTaP-IPV '
$ws.Range('A6').Value = 'DTaP-IPV '
$ws.Range('A7').Value = 'DTaP-IPV '
$ws.Range('A8').Value = 'DTaP-Hep B-IPV '
$ws.Range('A9').Value = 'DTaP-IP-HI '
$ws.Range('A10').Value = 'DTaP-IPV-HIB-HEPB '
$ws.Range('A11').Value = 'DTaP-IPV-HIB-HEPB '
$ws.Range('A12').Value = 'e-IPV '
$ws.Range('A13').Value = 'Hepatitis A Pediatric '
$ws.Range('A14').Value = 'Hepatitis A Pediatric '
$ws.Range('A15').Value = 'Hepatitis A-Hepatitis B 18 only '
$ws.Range('A16').Value = 'Hepatitis B  Pediatric/Adolescent'
$ws.Range('A17').Value = 'Hepatitis B  Pediatric/Adolescent'
$ws.Range('A18').Value = 'Hepatitis B  Pediatric/Adolescent'
$ws.Range('A19').Value = 'Hib '
$ws.Range('A20').Value = 'Hib '
$ws.Range('A21').Value = 'Hib '
$ws.Range('A22').Value = 'HPV - Human Papillomavirus 9-valent '
$ws.Range('A23').Value = 'MENB - Meningococcal Group B '
$ws.Range('A24').Value = 'MENB - Meningococcal Group B '
$ws.Range('A25').Value = 'Meningococcal Conjugate (Groups A, C, W and Y) '
$ws.Range('A26').Value = 'Meningococcal Conjugate (Groups A, C, Y and W-135) '
$ws.Range('A27').Value = 'Measles, Mumps and Rubella (MMR) '
$ws.Range('A28').Value = 'Measles, Mumps and Rubella (MMR) '
$ws.Range('A29').Value = 'MMR/Varicella '
$ws.Range('A30').Value = 'Pneumococcal 13-valent  (Pediatric)'
$ws.Range('A31').Value = 'Pneumococcal 15-valent '
$ws.Range('A33').Value = 'Rotavirus, Live, Oral, Pentavalent '
$ws.Range('A34').Value = 'Rotavirus, Live, Oral, Pentavalent '
$ws.Range('A35').Value = 'Rotavirus, Live, Oral, Oral '
$ws.Range('A36').Value = 'Tetanus and Diphtheria Toxoids '
$ws.Range('A37').Value = 'Tetanus and Diphtheria Toxoids '
$ws.Range('A38').Value = 'Tetanus and Diphtheria Toxoids '
$ws.Range('A39').Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws.Range('A40').Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws.Range('A41').Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws.Range('A42').Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws.Range('A43').Value = 'Varicella '

$ws = $wb.Worksheets.Item('Adult Vaccine ')
$ws.Range('A2').Value = 'Hepatitis A Adult '
$ws.Range('A3').Value = 'Hepatitis A Adult '
$ws.Range('A4').Value = 'Hepatitis A-Hepatitis B Adult '
$ws.Range('A5').Value = 'Hepatitis B Adult '
$ws.Range('A6').Value = 'Hepatitis B Adult '
$ws.Range('A7').Value = 'Hepatitis B Adult '
$ws.Range('A8').Value = 'Hepatitis B-Adult '
$ws.Range('A9').Value = 'HPV-Human Papillomavirus 9 Valent '
$ws.Range('A10').Value = 'Measles, Mumps,  Rubella '
$ws.Range('A11').Value = 'Measles, Mumps and Rubella (MMR) '
$ws.Range('A12').Value = 'Meningococcal Conjugate (Groups A, C, W and Y) '
$ws.Range('A13').Value = 'Meningococcal Conjugate (Groups A, C, Y and W-135) '
$ws.Range('A14').Value = 'MENB - Meningococcal Group B '
$ws.Range('A15').Value = 'MENB - Meningococcal Group B '
$ws.Range('A16').Value = 'Pneumococcal 15-valent '
$ws.Range('A17').Value = 'Pneumococcal 20-valent '
$ws.Range('A19').Value = 'Tetanus and Diphtheria Toxoids '
$ws.Range('A20').Value = 'Tetanus and Diphtheria Toxoids '
$ws.Range('A21').Value = 'Tetanus and Diphtheria Toxoids '
$ws.Range('A22').Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws.Range('A23').Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws.Range('A24').Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws.Range('A25').Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws.Range('A26').Value = 'Varicella '

$ws = $wb.Worksheets.Item('Pediatric Influenza Vaccine ')
$ws.Range('A2').Value = 'Influenza  (Age 6 months and older)'
$ws.Range('B2').Value = 'Fluzone Quadrivalent'
$ws.Range('A3').Value = 'Influenza  (Age 6 months and older)'
$ws.Range('B3').Value = 'Fluzone Quadrivalent'
$ws.Range('A4').Value = 'Influenza  (Age 6 months and older)'
$ws.Range('B4').Value = 'Fluzone Quadrivalent'
$ws.Range('A5').Value = 'Influenza  (Age 6 months and older)'
$ws.Range('B5').Value = 'Fluarix Quadrivalent'
$ws.Range('A6').Value = 'Influenza  (Age 6 months and older)'
$ws.Range('B6').Value = 'FluLaval Quadrivalent'
$ws.Range('A7').Value = 'Influenza  (Age 6 months and older)'
$ws.Range('A8').Value = 'Influenza  (Age 6 months and older)'
$ws.Range('A9').Value = 'Influenza  (Age 36 months and older)'
$ws.Range('A10').Value = 'Influenza  (Age 6 months and older)'
$ws.Range('A11').Value = 'Influenza  Live, Intranasal (Age 2-49 years)'
$ws.Range('B11').Value = 'FluMist Quadrivalent'

$ws = $wb.Worksheets.Item('Adult Influenza Vaccine ')
$ws.Range('A2').Value = 'Influenza  (Age 6 months and older)'
$ws.Range('B2').Value = 'Fluzone Quadrivalent'
$ws.Range('A3').Value = 'Influenza  (Age 6 months and older)'
$ws.Range('B3').Value = 'Fluzone Quadrivalent'
$ws.Range('A4').Value = 'Influenza  (Age 6 months and older)'
$ws.Range('B4').Value = 'Fluzone Quadrivalent'
$ws.Range('A5').Value = 'Influenza  (Age 6 months and older)'
$ws.Range('B5').Value = 'Fluarix Quadrivalent'
$ws.Range('A6').Value = 'Influenza  (Age 6 months and older)'
$ws.Range('B6').Value = 'FluLaval Quadrivalent'
$ws.Range('A7').Value = 'Influenza  (Age 6 months and older)'
$ws.Range('A8').Value = 'Influenza  (Age 6 months and older)'
$ws.Range('A9').Value = 'Influenza  (Age 36 months and older)'
$ws.Range('B9').Value = 'Afluria Quadrivalent'
$ws.Range('A10').Value = 'Influenza  (Age 6 months and older)'
$ws.Range('B10').Value = 'Afluria Quadrivalent'
$ws.Range('A11').Value = 'Influenza  Live, Intranasal (Age 2-49 years)'
$ws.Range('B11').Value = 'FluMist Quadrivalent'
